$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells are stored as text in the source sheet even though
# several look numeric (e.g. "208.44"). Excel auto-converts such strings to
# real numbers on assignment, so we flip those specific cells to Text format
# first to keep them as text, matching the original inline-string cell type.
$ws.Range("D2").Value = "27.570.77"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "1.596.82"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("E4").Value = "  +0.52%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.44"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("E6").Value = "  -3.45%  "
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("E8").Value = "  -4.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.252"
$ws.Range("E9").Value = "  -1.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0592"
$ws.Range("E10").Value = "  -3.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0862"
$ws.Range("E11").Value = "  -1.85%  "
$ws.Range("D12").Value = "1.823.87"
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("D13").Value = "1.610.99"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("E14").Value = "  -3.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.539"
$ws.Range("E15").Value = "  -3.95%  "
$ws.Range("E16").Value = "  -2.68%  "
$ws.Range("D17").Value = "27.573.20"
$ws.Range("E17").Value = "  -1.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "218.07"
$ws.Range("E18").Value = "  -4.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.40"
$ws.Range("E19").Value = "  -3.38%  "
$ws.Range("D20").Value = "0.0₃0695"
$ws.Range("E20").Value = "  -3.43%  "
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.21"
$ws.Range("E22").Value = "  -2.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.74"
$ws.Range("E23").Value = "  -3.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.00"
$ws.Range("E24").Value = "  -3.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.73"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("B26").Value = "BinanceUSD"
$ws.Range("C26").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.53%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.74"
$ws.Range("E27").Value = "  -1.97%  "
$ws.Range("E28").Value = "  -2.64%  "
$ws.Range("E29").Value = "  -3.78%  "
$ws.Range("E30").Value = "  -1.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0469"
$ws.Range("E31").Value = "  -2.40%  "
$ws.Range("E32").Value = "  -4.19%  "
$ws.Range("D33").Value = "1.369.40"
$ws.Range("E33").Value = "  -1.87%  "
$ws.Range("E34").Value = "  -4.43%  "
$ws.Range("E35").Value = "  -2.53%  "
$ws.Range("E36").Value = "  -5.77%  "
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("E38").Value = "  -2.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.540"
$ws.Range("E39").Value = "  -2.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.814"
$ws.Range("E40").Value = "  -4.32%  "
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.979"
$ws.Range("E42").Value = "  -2.77%  "
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.78"
$ws.Range("E44").Value = "  -3.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.15"
$ws.Range("E45").Value = "  -2.36%  "
$ws.Range("D46").Value = "1.734.45"
$ws.Range("E46").Value = "  -1.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.10"
$ws.Range("E47").Value = "  -2.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.00"
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").Value = "0.0₇0999"
$ws.Range("E49").Value = "  -0.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0972"
$ws.Range("E50").Value = "  -4.16%  "
$ws.Range("E51").Value = "  -0.89%  "
